$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string values that changed
$ws.Range("B2").Value = "ff7Zo5"
$ws.Range("U2").Value = "-"
$ws.Range("V2").Value = "-"

# V2's style switches from right-aligned (matching numeric cells) to
# left-aligned (matching the "oip" column U2) since it now shares the
# same text value/style as U2.
$ws.Range("V2").HorizontalAlignment = -4131
